# 13.Списъци-увод.pptx — "Improved the slides (formatting)"
#
# 1. Remove all speaker/reviewer comments (comment1..comment12) from every slide.
# 2. Slide 13 ("Какво научихме този час?"):
#    - Content placeholder: drop the lnSpcReduction on normAutofit (text no
#      longer needs manual shrinking), retitle "Lists" -> "Списъкът" (and mark
#      that run as Bulgarian), and collapse the two-line "Изпечатване на /
#      елементите на списък:" into the single line "Отпечатване на списък:".
#    - Reposition/resize the two pictures and the three code-sample callouts
#      lower/wider on the slide.

$p = $ppt.ActivePresentation

# --- 1. Strip all the reviewer comments across the whole deck -------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    while ($s.Comments.Count -gt 0) {
        $s.Comments.Item(1).Delete()
    }
}

# --- 2. Slide 13 content + layout tweaks -----------------------------------
$slide = $p.Slides.Item(13)

# Content Placeholder 4 (id=5) — the bullet list.
$body = $slide.Shapes.Item(1)
$body.TextFrame.AutoSize = 2   # -> <a:normAutofit/> (no lnSpcReduction)

$tr = $body.TextFrame.TextRange

$listsRun = $tr.Find("Lists", 0, $false)
$listsRun.LanguageID = "bg-BG"
$listsRun.Text = "Списъкът"

$printStart = $tr.Find("Изпечатване", 0, $false)
$printRange = $tr.Characters($printStart.Start, 36)
$printRange.Text = "Отпечатване на списък:"

# Picture 9 (id=10)
$pic9 = $slide.Shapes.Item(3)
$pic9.Left = 682.071533203125
$pic9.Top = 162.8810272216797

# Picture 10 (id=11)
$pic10 = $slide.Shapes.Item(4)
$pic10.Left = 760.9922485351562
$pic10.Top = 221.8678741455078

# Text Placeholder 5 (id=12) — first code callout box
$callout1 = $slide.Shapes.Item(5)
$callout1.Top = 312.0
$callout1.Width = 611.5689086914062

# Text Placeholder 5 (id=13) — second code callout box
$callout2 = $slide.Shapes.Item(6)
$callout2.Left = 362.63922119140625
$callout2.Top = 406.40655517578125

# Text Placeholder 5 (id=14) — third code callout box
$callout3 = $slide.Shapes.Item(7)
$callout3.Left = 362.88671875
$callout3.Top = 465.4630126953125
$callout3.Width = 563.9998779296875
